# Applies the cryptos price/volume refresh described in the commit diff.
# Every D/E cell in this sheet is a plain text cell (t="inlineStr" in the
# source OOXML) -- e.g. "67.505.27" or "  -1.42%  " -- never a real number.
# Excel/COM auto-coerces any Value string that *parses* as a number (like
# "600.30" or "27.86") into a Double, which both changes the stored cell type
# and can silently drop a trailing zero (600.30 -> 600.3). To keep those cells
# textual -- matching the target content exactly -- we briefly force a text
# NumberFormat before assigning, then clear the format again so the cell keeps
# using the sheet default style (no explicit s="..." attribute), exactly like
# every other untouched data cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.617.43"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "2.680.66"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "2.680.08"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.359"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.86"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "3.168.03"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "67.510.55"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "2.669.48"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.83"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "558.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.94"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.58"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.37"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.375"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.95"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "0.0₆0298"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.59"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.83"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.22%  "
